# Apply cryptos.xlsx price/volume update (GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.345.85"
$ws.Range("E2").Value = "  -2.14%  "
$ws.Range("D3").Value = "2.880.79"
$ws.Range("E3").Value = "  -1.87%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.17%  "
$ws.Range("D5").Value = "'567.15"
$ws.Range("E5").Value = "  -4.36%  "
$ws.Range("D6").Value = "'143.39"
$ws.Range("E6").Value = "  -1.92%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").Value = "'0.505"
$ws.Range("E8").Value = "  -0.11%  "
$ws.Range("D9").Value = "2.878.86"
$ws.Range("E9").Value = "  -1.89%  "
$ws.Range("D10").Value = "'6.89"
$ws.Range("E10").Value = "  -5.58%  "
$ws.Range("D11").Value = "'0.147"
$ws.Range("E11").Value = "  -1.84%  "
$ws.Range("E12").Value = "  -1.99%  "
$ws.Range("D13").Value = "'0.0000234"
$ws.Range("E13").Value = "  -0.80%  "
$ws.Range("D14").Value = "'31.90"
$ws.Range("E14").Value = "  -2.06%  "
$ws.Range("E15").Value = "  -0.44%  "
$ws.Range("D16").Value = "3.349.67"
$ws.Range("E16").Value = "  -2.16%  "
$ws.Range("D17").Value = "61.280.44"
$ws.Range("E17").Value = "  -2.23%  "
$ws.Range("D18").Value = "'6.57"
$ws.Range("E18").Value = "  -1.36%  "
$ws.Range("D19").Value = "2.882.76"
$ws.Range("E19").Value = "  -1.90%  "
$ws.Range("D20").Value = "'430.44"
$ws.Range("E20").Value = "  -1.75%  "
$ws.Range("D21").Value = "'13.04"
$ws.Range("E21").Value = "  -2.14%  "
$ws.Range("D22").Value = "'0.654"
$ws.Range("E22").Value = "  -1.36%  "
$ws.Range("D23").Value = "'6.83"
$ws.Range("E23").Value = "  -2.60%  "
$ws.Range("D24").Value = "'79.06"
$ws.Range("E24").Value = "  -2.25%  "
$ws.Range("D25").Value = "'11.78"
$ws.Range("E25").Value = "  +0.79%  "
$ws.Range("E26").Value = "  +0.14%  "
$ws.Range("D27").Value = "'10.03"
$ws.Range("E27").Value = "  -8.86%  "
$ws.Range("D28").Value = "'2.01"
$ws.Range("E28").Value = "  -4.96%  "
$ws.Range("D29").Value = "'0.0000104"
$ws.Range("E29").Value = "  +3.03%  "
$ws.Range("D30").Value = "'6.98"
$ws.Range("E30").Value = "  -2.65%  "
$ws.Range("E31").Value = "  -4.37%  "
$ws.Range("E32").Value = "  -7.88%  "
$ws.Range("D33").Value = "'1.00"
$ws.Range("E33").Value = "  -0.11%  "
$ws.Range("D34").Value = "'0.106"
$ws.Range("E34").Value = "  -2.22%  "
$ws.Range("D35").Value = "'25.50"
$ws.Range("E35").Value = "  -3.01%  "
$ws.Range("D36").Value = "'0.956"
$ws.Range("E36").Value = "  -3.60%  "
$ws.Range("D37").Value = "'5.38"
$ws.Range("E37").Value = "  -2.93%  "
$ws.Range("D38").Value = "'48.81"
$ws.Range("E38").Value = "  -1.61%  "
$ws.Range("B39").Value = "dogwifhat"
$ws.Range("C39").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D39").Value = "'2.83"
$ws.Range("E39").Value = "  -8.22%  "
$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").Value = "'1.94"
$ws.Range("E40").Value = "  -3.82%  "
$ws.Range("D41").Value = "'8.21"
$ws.Range("E41").Value = "  -2.74%  "
$ws.Range("E42").Value = "  -2.86%  "
$ws.Range("D43").Value = "'39.14"
$ws.Range("E43").Value = "  +0.42%  "
$ws.Range("D44").Value = "'0.268"
$ws.Range("E44").Value = "  -3.89%  "
$ws.Range("D45").Value = "2.695.55"
$ws.Range("E45").Value = "  -0.02%  "
$ws.Range("D46").Value = "'133.87"
$ws.Range("E46").Value = "  -0.32%  "
$ws.Range("D47").Value = "'0.0334"
$ws.Range("E47").Value = "  -0.29%  "
$ws.Range("E48").Value = "  -0.05%  "
$ws.Range("D49").Value = "'339.29"
$ws.Range("E49").Value = "  -5.87%  "
$ws.Range("D50").Value = "'0.103"
$ws.Range("E50").Value = "  -1.25%  "
$ws.Range("D51").Value = "'21.51"
$ws.Range("E51").Value = "  -5.00%  "

# Clear the quote-prefix flag Excel applies when a text value looks numeric,
# so these cells keep the workbook default (unstyled) formatting.
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D18").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Style = "Normal"
